# The "Mapping CSV2openEHR" sheet holds the CSV -> openEHR FLAT path mapping.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping CSV2openEHR")
$ws.Activate()

# Row 5 (administratives_geschlecht_value) previously had no FLAT-Path chosen
# from the dropdown (data validation list in column B). Pick the matching
# path for the context start time, same as already used in row 6 of
# FLAT_Paths ("bericht/context/start_time").
$ws.Range("B5").Value = "bericht/context/start_time"

# Leave the selection on B13, reflecting where editing last left off.
$ws.Range("B13").Select()

$wb.Save()
